$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 21 for the Wenaha River / WEN site (shifts old rows 21-24 down to 22-25)
$ws.Rows.Item(21).Insert()

# Populate the appended Tucannon River rows (26-30) first so their new shared
# strings are registered ahead of the Wenaha River row's strings.
$tucannonSites = "LTR", "MTR", "UTR", "TFH", "TPJ"
$row = 26
foreach ($site in $tucannonSites) {
    $ws.Range("A" + $row).Value = $site
    $row = $row + 1
}
$row = 26
foreach ($site in $tucannonSites) {
    $ws.Range("B" + $row).Value = "Snake River Coho Salmon"
    $ws.Range("C" + $row).Value = "Lower Snake"
    $ws.Range("D" + $row).Value = "SNTUC-c"
    $ws.Range("E" + $row).Value = "Tucannon River"
    $row = $row + 1
}

$ws.Range("A21").Value = "WEN"
$ws.Range("B21").Value = "Snake River Coho Salmon"
$ws.Range("C21").Value = "Grande Ronde River"
$ws.Range("D21").Value = "GRWEN-c"
$ws.Range("E21").Value = "Wenaha River"
